# Updates crypto price ("D" column) and 1h-volume percentage ("E" column)
# values on the active worksheet to match the latest scrape.
# Values are stored as plain text (matching the existing inline-string
# cells), so the cells' NumberFormat is forced to "text" before writing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @(newPrice, newVolumePercent)
# An empty string means "leave the existing value untouched".
$updates = @{
    2  = @("260.46", "1.44%")
    3  = @("27.23", "2.31%")
    4  = @("", "0.65%")
    5  = @("0.06172", "4.19%")
    6  = @("6.662", "0.86%")
    7  = @("0.8521", "-0.41%")
    8  = @("0.9170", "0.48%")
    9  = @("0.1412", "")
    10 = @("0.04825", "9.43%")
    11 = @("0.07081", "1.07%")
    12 = @("0.03110", "3.24%")
    13 = @("0.09043", "-0.63%")
    14 = @("0.001535", "0.04%")
    15 = @("0.0006181", "-94.04%")
    16 = @("0.006044", "-1.22%")
    18 = @("3.155", "0.77%")
    19 = @("", "1.39%")
    20 = @("", "-0.34%")
    21 = @("0.1300", "0.32%")
    22 = @("4.090", "6.16%")
    23 = @("0.04254", "1.06%")
    24 = @("", "0.12%")
    25 = @("0.003802", "")
    26 = @("", "0.11%")
    27 = @("", "-8.08%")
    40 = @("0.03873", "1.91%")
    41 = @("0.1114", "1.31%")
    42 = @("0.004090", "-34.05%")
    43 = @("", "8.86%")
    44 = @("0.002204", "0.30%")
    45 = @("0.00005166", "1.63%")
    46 = @("", "0.11%")
    47 = @("", "8.16%")
    48 = @("0.1229", "-48.95%")
    49 = @("", "0.11%")
    50 = @("", "0.11%")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($priceVal -ne "") {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $priceVal
    }
    if ($volVal -ne "") {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $volVal
    }
}
